$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 37-50 with new chat log entries (ubike/bus info)
$ws.Range("D37").Value = "2022-05-08 20:54:41.028000"

$ws.Range("B38").Value = "bus 你好"
$ws.Range("C38").Value = "桃園公車中查無此資料"
$ws.Range("D38").Value = "2022-05-08 20:58:06.508000"

$ws.Range("B39").Value = "ubike 哈哈"
$ws.Range("C39").Value = "桃園ubike中查無此資料"
$ws.Range("D39").Value = "2022-05-08 20:58:22.101000"

$ws.Range("B40").Value = "鋼筆多少"
$ws.Range("C40").Value = "鋼筆 30元"
$ws.Range("D40").Value = "2022-05-08 21:00:10.521000"

$ws.Range("B41").Value = "怎麼購買"
$ws.Range("C41").Value = "nontextreply"
$ws.Range("D41").Value = "2022-05-08 21:00:28.830000"

$ws.Range("B42").Value = "bus 健行科技大學"
$ws.Range("C42").Value = "桃園公車中查無此資料"
$ws.Range("D42").Value = "2022-05-08 21:01:03.663000"

$ws.Range("B43").Value = "ubike 健行科技大學"
$ws.Range("C43").Value = "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:34`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1"
$ws.Range("D43").Value = "2022-05-08 21:01:22.431000"

$ws.Range("B44").Value = "Bus 171-FS"
$ws.Range("D44").Value = "2022-05-08 21:12:24.115000"

$ws.Range("B45").Value = "bus 171-FS "
$ws.Range("C45").Value = "車輛:171-FS`n業者代號:45`nGPS車速:0.0`nGPS時間:2022-05-07 22:15:26`n路線方向(1:去程,2:回程):1"
$ws.Range("D45").Value = "2022-05-08 21:12:45.175000"

$ws.Range("B46").Value = "位置在哪"
$ws.Range("D46").Value = "2022-05-08 21:22:06.813000"

$ws.Range("B47").Value = "位置"
$ws.Range("C47").Value = "nontextreply"
$ws.Range("D47").Value = "2022-05-08 21:22:44.916000"

$ws.Range("B48").Value = "位置呢"
$ws.Range("D48").Value = "2022-05-08 21:23:01.917000"

$ws.Range("B49").Value = "地址是？"
$ws.Range("C49").Value = "261宜蘭縣頭城鎮港口路92-1號"
$ws.Range("D49").Value = "2022-05-08 21:23:29.902000"

$ws.Range("B50").Value = "ubike 健行科技大學"
$ws.Range("C50").Value = "中文場站名稱:健行科技大學`n場站總停車格:66`n場站目前車輛數:34`n地址:健行路229號(商學大樓後人行道)`n場站是否暫停營運1"
$ws.Range("D50").Value = "2022-05-08 21:25:29.317000"

# Remove trailing rows 51-56 that no longer exist in the updated log
$ws.Rows("51:56").Delete()

# Re-fit the height of rows whose new content spans multiple lines
$ws.Rows("43:45").AutoFit()
$ws.Rows("50:50").AutoFit()

# Best-fit the column widths for the (now wider/narrower) data
$ws.Columns.Item(1).ColumnWidth = 36.58515625
$ws.Columns.Item(2).ColumnWidth = 14.440625
$ws.Columns.Item(3).ColumnWidth = 36.87031250
$ws.Columns.Item(4).ColumnWidth = 26.58515625

# Restore the view to what the author left it at while reviewing the new rows
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D31").Select()
